$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows (C4:N13) ---------------------------------------------
$rows = @(
    @{ id=1;  colab=1; nome="Marlon Passeri";          mes="2025-06"; horasMes=10;   gp=656;  horasTrab=10;   prop=1 },
    @{ id=2;  colab=2; nome="Leonardo";                 mes="2025-06"; horasMes=14;   gp=1796; horasTrab=10;   prop=0.7142857142857143 },
    @{ id=3;  colab=2; nome="Leonardo";                 mes="2025-06"; horasMes=14;   gp=1003; horasTrab=4;    prop=0.2857142857142857 },
    @{ id=4;  colab=3; nome="Fel";                      mes="2025-06"; horasMes=8;    gp=2516; horasTrab=8;    prop=1 },
    @{ id=5;  colab=4; nome="Marlon Soares Passeri";    mes="2025-06"; horasMes=9;    gp=1881; horasTrab=9;    prop=1 },
    @{ id=6;  colab=5; nome="Guilherme Fernandes";      mes="2025-06"; horasMes=32.5; gp=656;  horasTrab=4;    prop=0.1230769230769231 },
    @{ id=7;  colab=5; nome="Guilherme Fernandes";      mes="2025-06"; horasMes=32.5; gp=669;  horasTrab=4;    prop=0.1230769230769231 },
    @{ id=8;  colab=5; nome="Guilherme Fernandes";      mes="2025-06"; horasMes=32.5; gp=1796; horasTrab=14.5; prop=0.4461538461538462 },
    @{ id=9;  colab=5; nome="Guilherme Fernandes";      mes="2025-06"; horasMes=32.5; gp=3040; horasTrab=1;    prop=0.03076923076923077 },
    @{ id=10; colab=5; nome="Guilherme Fernandes";      mes="2025-06"; horasMes=32.5; gp=2139; horasTrab=9;    prop=0.2769230769230769 }
)

$r = 4
foreach ($row in $rows) {
    $ws.Cells.Item($r, 3).Value  = $row.id        # C - ID
    $ws.Cells.Item($r, 4).Value  = $row.colab      # D - Id Colaborador
    $ws.Cells.Item($r, 5).Value  = $row.nome       # E - Nome Colaborador
    $ws.Cells.Item($r, 6).Value  = $row.mes        # F - Mes
    $ws.Cells.Item($r, 7).Value  = $row.horasMes   # G - Horas Mes
    $ws.Cells.Item($r, 8).Value  = $row.gp         # H - GP
    $ws.Cells.Item($r, 9).Value  = $row.horasTrab  # I - Horas Trabalhadas
    $ws.Cells.Item($r, 10).Value = $row.prop       # J - Proporcao de Hora
    $r = $r + 1
}

# --- Borders: bottom border on every data row, left edge on col C, right edge on col N.
# (Bottom is applied first, for the whole C:N width of every row, so that the
# later single-side additions on column C/N extend that already-committed
# border instead of going through an unrelated intermediate combination.)
for ($row = 4; $row -le 13; $row++) {
    $ws.Range("C" + $row + ":N" + $row).Borders.Item(9).LineStyle = 1
}
$ws.Range("C4:C13").Borders.Item(7).LineStyle = 1
$ws.Range("N4:N13").Borders.Item(10).LineStyle = 1

# --- Number formats (applied after borders so each number-format xf already
# carries the right borderId instead of allocating a throw-away border-less one)
$ws.Range("G4:G13").NumberFormat = "0.00"
$ws.Range("I4:I13").NumberFormat = "0.00"
$ws.Range("J4:J13").NumberFormat = "0.00%"

Write-Output "applied edits"
